$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Title paragraph: "Service Learning" / " Diary " (split across two runs
#    around a proofErr gramStart/gramEnd pair) -> single run, proofErr gone.
#    The gramStart marker sits at the very start of the paragraph (position
#    0 of the match), so a plain Find/Replace over the paragraph text won't
#    absorb it - insert a throw-away run just before the paragraph first so
#    the marker ends up strictly *inside* the matched/replaced range, then
#    let Find.Execute fold everything (dummy run + proofErr + the two text
#    runs) into one clean run.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(2)
$titlePara.Range.InsertBefore("ZZDUMMYZZ")
$d.Content.Find.Execute("ZZDUMMYZZService Learning Diary", $false, $false, $false, $false, $false, $true, 1, $false, "Service Learning Diary", 2)

# ---------------------------------------------------------------------------
# 2) "... about meetings (with service learning lecturers, ingroup and with
#    stakeholders)" - proofErr pair is preceded by real text in the same
#    run-group, so a normal whole-phrase Find/Replace merges the runs and
#    drops both proofErr markers.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(" about meetings (with service learning lecturers, ingroup and with stakeholders)", $false, $false, $false, $false, $false, $true, 1, $false, " about meetings (with service learning lecturers, ingroup and with stakeholders)", 2)

# ---------------------------------------------------------------------------
# 3) "... during literature research/practical research e.g. at site
#    (observations, talks etc.)"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(" during literature research/practical research e.g. at site (observations, talks etc.)", $false, $false, $false, $false, $false, $true, 1, $false, " during literature research/practical research e.g. at site (observations, talks etc.)", 2)

# ---------------------------------------------------------------------------
# 4) "What are key insights I gained during the meetings (with service
#    learning lecturers, ingroup and with stakeholders)?"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("What are key insights I gained during the meetings (with service learning lecturers, ingroup and with stakeholders)?", $false, $false, $false, $false, $false, $true, 1, $false, "What are key insights I gained during the meetings (with service learning lecturers, ingroup and with stakeholders)?", 2)

# ---------------------------------------------------------------------------
# 5) "What different hypotheses did I develop during the interaction with
#    service learning lecturers, ingroup and with stakeholders?"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("What different hypotheses did I develop during the interaction with service learning lecturers, ingroup and with stakeholders?", $false, $false, $false, $false, $false, $true, 1, $false, "What different hypotheses did I develop during the interaction with service learning lecturers, ingroup and with stakeholders?", 2)

# ---------------------------------------------------------------------------
# 6) "What did I learn during the stakeholder analysis, the project problem
#    description, etc\u2026 " - spellStart/gramStart/spellEnd/gramEnd all
#    collapse away once the whole phrase (through the trailing space) is
#    matched/replaced in one go.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("What did I learn during the stakeholder analysis, the project problem description, etc… ", $false, $false, $false, $false, $false, $true, 1, $false, "What did I learn during the stakeholder analysis, the project problem description, etc… ", 2)

# ---------------------------------------------------------------------------
# 7) Header: "GMIT Civic Engagement" -> "Civic Engagement"
# ---------------------------------------------------------------------------
$section = $d.Sections.Item(1)
$header = $section.Headers.Item(1)
$header.Range.Find.Execute("GMIT Civic Engagement", $false, $false, $false, $false, $false, $true, 1, $false, "Civic Engagement", 2)
